# Weekly update: insert a new price record for Alcachofa / Madrigal at row 36,
# pushing the existing rows 36-45 down to 37-46 (one extra row of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 36; this shifts rows 36-45 down to 37-46,
# carrying their existing formatting and values with them.
$ws.Rows.Item(36).Insert()

# Fill the newly inserted row 36 with the new weekly data point.
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 45204
$ws.Cells.Item(36, 4).NumberFormat = $ws.Cells.Item(37, 4).NumberFormat
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112013
$ws.Cells.Item(36, 7).Value = "Alcachofa"
$ws.Cells.Item(36, 8).Value = "Madrigal"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 140
$ws.Cells.Item(36, 11).Value = 15000
$ws.Cells.Item(36, 12).Value = 16000
$ws.Cells.Item(36, 13).Value = 15643
$ws.Cells.Item(36, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(36, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(36, 16).Value = 391
$ws.Cells.Item(36, 17).Value = 40
$ws.Cells.Item(36, 18).Value = "Hortaliza"
